# Rename the 50 summary sheets to reflect the new spatial-unit IDs
# produced by the rerun of the LU d2c FeatEng pipeline for FR cities.
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ29416723","summ29604762","summ29788548","summ29977503","summ30171389",
    "summ30368842","summ30561157","summ30755071","summ30971902","summ31240680",
    "summ31505169","summ31775463","summ32056301","summ32322137","summ32590861",
    "summ32871889","summ33140304","summ33407530","summ33673712","summ33939753",
    "summ34203609","summ34489358","summ34772652","summ35038623","summ35313202",
    "summ35616160","summ35933469","summ36238162","summ36523407","summ36773772",
    "summ37074961","summ37391755","summ37706103","summ38022916","summ38305047",
    "summ38594710","summ38971230","summ39225007","summ39481261","summ39741514",
    "summ40005172","summ40291705","summ40573454","summ40875080","summ41138784",
    "summ41440107","summ41688620","summ41972177","summ42255601","summ42530850"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
